# Actualizacion automatica de scrims_actualizado.xlsx (2025-07-23 23:43:20)
# Appends new scrim result rows to three sheets: "Crystal Arcade",
# "New Horizons" and "Hot Potato".

$wb = $excel.ActiveWorkbook

function Add-ScrimRow($worksheet, $newRow, $formatSourceRow, $values) {
    # Copy formatting (fills/borders/fonts -> cell styles) from an existing
    # row that already carries the right "Equipo 1"/"Equipo 2" styling, so
    # the style table itself stays untouched.
    $srcRange = $worksheet.Range("A" + $formatSourceRow + ":N" + $formatSourceRow)
    $dstRange = $worksheet.Range("A" + $newRow + ":N" + $newRow)
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial(-4122) | Out-Null

    for ($col = 1; $col -le 14; $col++) {
        $worksheet.Cells.Item($newRow, $col).Value = $values[$col - 1]
    }
}

# ---------------------------------------------------------------------
# Sheet "Crystal Arcade": add row 29 (Equipo 1 style, like row 4)
# ---------------------------------------------------------------------
$wsCrystal = $wb.Worksheets.Item("Crystal Arcade")
Add-ScrimRow $wsCrystal 29 4 @(
    "GRAY",
    "FANG",
    "EMZ",
    "ALLI",
    "SANDY",
    "GUS",
    "Equipo 1",
    "TRB|Zeus 解開",
    "TRB|Lxffy",
    "TRB|R B M",
    "NHG|Xemp",
    "KCP|Fade",
    "KCP|Tyrant",
    "20250723T214219.000Z"
)

# ---------------------------------------------------------------------
# Sheet "New Horizons": add rows 34-39
# ---------------------------------------------------------------------
$wsHorizons = $wb.Worksheets.Item("New Horizons")

Add-ScrimRow $wsHorizons 34 4 @(
    "ANGELO",
    "JUJU",
    "KAZE",
    "R-T",
    "DARRYL",
    "SQUEAK",
    "Equipo 1",
    "KCP|Fade",
    "KCP|Tyrant",
    "NHG|Xemp",
    "TRB|Zeus 解開",
    "TRB|R B M",
    "TRB|Lxffy",
    "20250723T213441.000Z"
)

Add-ScrimRow $wsHorizons 35 4 @(
    "R-T",
    "SQUEAK",
    "DARRYL",
    "ANGELO",
    "KAZE",
    "JUJU",
    "Equipo 1",
    "TRB|Zeus 解開",
    "TRB|Lxffy",
    "TRB|R B M",
    "KCP|Fade",
    "NHG|Xemp",
    "KCP|Tyrant",
    "20250723T213259.000Z"
)

Add-ScrimRow $wsHorizons 36 8 @(
    "R-T",
    "SQUEAK",
    "COLT",
    "KAZE",
    "ANGELO",
    "JUJU",
    "Equipo 2",
    "TRB|Zeus 解開",
    "TRB|Lxffy",
    "TRB|R B M",
    "NHG|Xemp",
    "KCP|Fade",
    "KCP|Tyrant",
    "20250723T213012.000Z"
)

Add-ScrimRow $wsHorizons 37 8 @(
    "R-T",
    "SQUEAK",
    "COLT",
    "KAZE",
    "ANGELO",
    "JUJU",
    "Equipo 2",
    "TRB|Zeus 解開",
    "TRB|Lxffy",
    "TRB|R B M",
    "NHG|Xemp",
    "KCP|Fade",
    "KCP|Tyrant",
    "20250723T212737.000Z"
)

Add-ScrimRow $wsHorizons 38 4 @(
    "GROM",
    "CORDELIUS",
    "GENE",
    "TICK",
    "BUSTER",
    "MANDY",
    "Equipo 1",
    "TRB|Zeus 解開",
    "TRB|Lxffy",
    "TRB|R B M",
    "NHG|Xemp",
    "KCP|Fade",
    "KCP|Tyrant",
    "20250723T212215.000Z"
)

Add-ScrimRow $wsHorizons 39 4 @(
    "GROM",
    "CORDELIUS",
    "GENE",
    "TICK",
    "BUSTER",
    "MANDY",
    "Equipo 1",
    "TRB|Zeus 解開",
    "TRB|Lxffy",
    "TRB|R B M",
    "NHG|Xemp",
    "KCP|Fade",
    "KCP|Tyrant",
    "20250723T212017.000Z"
)

# ---------------------------------------------------------------------
# Sheet "Hot Potato": add row 45 (Equipo 2 style, like row 4)
# ---------------------------------------------------------------------
$wsHotPotato = $wb.Worksheets.Item("Hot Potato")
Add-ScrimRow $wsHotPotato 45 4 @(
    "OTIS",
    "RUFFS",
    "HANK",
    "KAZE",
    "CROW",
    "GRIFF",
    "Equipo 2",
    "TRB|Zeus 解開",
    "TRB|Lxffy",
    "TRB|R B M",
    "NHG|Xemp",
    "KCP|Fade",
    "KCP|Tyrant",
    "20250723T211205.000Z"
)
